$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "1942 (never released in Japan)"
$ws.Range("H2").Value = "Palette 0x03 (2 colors) on purpose. You can print scores every 4 levels or at game over"

# Row 39 - "Hello Kitty no Magical museum"
$ws.Range("H39").Value = "Palette 0xE4. You can print about 45 cute images after completing puzzle games."

# Row 57 - "McDonald's Monogatari : Honobono Tenchou Ikusei Game"
$ws.Range("H57").Value = "Palette 0x30 (2 colors). You can print the result from the cashier minigame."

# Row 63 - "Nakayoshi Cooking Series 1 - Oishii Cake-ya-san"
$ws.Range("H63").Value = "Palette 0x07 (2 colors), you can print recipes.."

# Row 64 - "Nakayoshi Cooking Series 2 - Oishii Panya-san"
$ws.Range("H64").Value = "Palette 0x07 (2 colors), you can print recipes."

# Row 102 - "Tales of Phantasia: Nakiri's Dungeon"
$ws.Range("H102").Value = "Palette 0x00, acts as 0xE4 (documented in pandocs). You can print images of your team characters in diffrent costumes. Printing protocol spams 0F commands but works."
